$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (B2): currency wording MPLN/PLN -> MEUR/EUR ---
$b2 = $ws.Range("B2")
$b2.Value = "ETS CO2 Price Projection [MEUR/kt = EUR/kg]"

# Re-apply the original rich-text look (white Arial 11, "2" as a subscript)
# to the new title text so the cell keeps its banner styling.
$b2.Characters(1, 6).Font.Name = "Arial"
$b2.Characters(1, 6).Font.Size = 11
$b2.Characters(1, 6).Font.Color = 16777215

$b2.Characters(7, 1).Font.Name = "Arial"
$b2.Characters(7, 1).Font.Size = 11
$b2.Characters(7, 1).Font.Color = 16777215
$b2.Characters(7, 1).Font.Subscript = $true

$b2.Characters(8, 37).Font.Name = "Arial"
$b2.Characters(8, 37).Font.Size = 11
$b2.Characters(8, 37).Font.Color = 16777215

# --- Region scope (D5): Poland-only -> all regions ---
$ws.Range("D5").Value = "AllRegions"

# --- Remove the blank separator row (old row 7) so the table closes up ---
$ws.Rows(7).Delete()

# --- Update the active selection to match the saved view ---
$ws.Range("K6").Select()
